$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Row 7: "Experimental" -> set Value (column B) to true (leading apostrophe forces text, not boolean)
$ws.Range("B7").Value = "'true"

# Row 8: "Date" -> update the Value (column B) to new timestamp
$ws.Range("B8").Value = "2025-01-28T15:58:19+00:00"
